$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.0271363451273
$ws.Cells.Item(2, 4).Value = 1.031593715891273
$ws.Cells.Item(2, 5).Value = 1.050376297517652
$ws.Cells.Item(2, 6).Value = 1.05525434559609
$ws.Cells.Item(2, 9).Value = 1.032316750418348
$ws.Cells.Item(2, 10).Value = 1.032295897766756
$ws.Cells.Item(2, 11).Value = 1.034401384282966
$ws.Cells.Item(2, 12).Value = 1.053130746599599
$ws.Cells.Item(2, 13).Value = 1.057995308916699
$ws.Cells.Item(2, 14).Value = 1.014719588740343

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.028303656739852
$ws.Cells.Item(3, 4).Value = 1.032460587905173
$ws.Cells.Item(3, 5).Value = 1.051563489139146
$ws.Cells.Item(3, 6).Value = 1.056493664502133
$ws.Cells.Item(3, 9).Value = 1.032542132403142
$ws.Cells.Item(3, 10).Value = 1.033102146747051
$ws.Cells.Item(3, 11).Value = 1.035076696333033
$ws.Cells.Item(3, 12).Value = 1.054129406533073
$ws.Cells.Item(3, 13).Value = 1.059046946189337
$ws.Cells.Item(3, 14).Value = 1.014992280135225

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.029058574355807
$ws.Cells.Item(4, 4).Value = 1.033020964352403
$ws.Cells.Item(4, 5).Value = 1.052332222379775
$ws.Cells.Item(4, 6).Value = 1.057295953330431
$ws.Cells.Item(4, 9).Value = 1.03268641768842
$ws.Cells.Item(4, 10).Value = 1.033622936928307
$ws.Cells.Item(4, 11).Value = 1.03551248614413
$ws.Cells.Item(4, 12).Value = 1.054775574297346
$ws.Cells.Item(4, 13).Value = 1.05972722942646
$ws.Cells.Item(4, 14).Value = 1.015168276883703

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.029375845113489
$ws.Cells.Item(5, 4).Value = 1.033256415940301
$ws.Cells.Item(5, 5).Value = 1.052655528588349
$ws.Cells.Item(5, 6).Value = 1.057633324741576
$ws.Cells.Item(5, 9).Value = 1.032746703852003
$ws.Cells.Item(5, 10).Value = 1.03384166093166
$ws.Cells.Item(5, 11).Value = 1.035695409713108
$ws.Cells.Item(5, 12).Value = 1.05504721669227
$ws.Cells.Item(5, 13).Value = 1.060013174688135
$ws.Cells.Item(5, 14).Value = 1.015242157718477

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.029429110656874
$ws.Cells.Item(6, 4).Value = 1.033295941648431
$ws.Cells.Item(6, 5).Value = 1.052709820827742
$ws.Cells.Item(6, 6).Value = 1.057689976157232
$ws.Cells.Item(6, 9).Value = 1.032756804391566
$ws.Cells.Item(6, 10).Value = 1.033878373003111
$ws.Cells.Item(6, 11).Value = 1.03572610683751
$ws.Cells.Item(6, 12).Value = 1.055092826295289
$ws.Cells.Item(6, 13).Value = 1.060061183504497
$ws.Cells.Item(6, 14).Value = 1.015254556296919

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.029062814120558
$ws.Cells.Item(7, 4).Value = 1.033024110980938
$ws.Cells.Item(7, 5).Value = 1.052336541897794
$ws.Cells.Item(7, 6).Value = 1.05730046095152
$ws.Cells.Item(7, 9).Value = 1.032687224694119
$ws.Cells.Item(7, 10).Value = 1.033625860377496
$ws.Cells.Item(7, 11).Value = 1.035514931487026
$ws.Cells.Item(7, 12).Value = 1.054779204020585
$ws.Cells.Item(7, 13).Value = 1.059731050419541
$ws.Cells.Item(7, 14).Value = 1.015169264507922

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.027530929747106
$ws.Cells.Item(8, 4).Value = 1.031886793029
$ws.Cells.Item(8, 5).Value = 1.050777403255422
$ws.Cells.Item(8, 6).Value = 1.055673104163698
$ws.Cells.Item(8, 9).Value = 1.032393240752521
$ws.Cells.Item(8, 10).Value = 1.032568561694526
$ws.Cells.Item(8, 11).Value = 1.034629854283923
$ws.Cells.Item(8, 12).Value = 1.053468255701445
$ws.Cells.Item(8, 13).Value = 1.058350756049397
$ws.Cells.Item(8, 14).Value = 1.014811839958557

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02482832464815
$ws.Cells.Item(9, 4).Value = 1.029878477641648
$ws.Cells.Item(9, 5).Value = 1.048034087184047
$ws.Cells.Item(9, 6).Value = 1.052808237729338
$ws.Cells.Item(9, 9).Value = 1.031863309861523
$ws.Cells.Item(9, 10).Value = 1.030698480075771
$ws.Cells.Item(9, 11).Value = 1.033061154781701
$ws.Cells.Item(9, 12).Value = 1.051157901566162
$ws.Cells.Item(9, 13).Value = 1.055916942441552
$ws.Cells.Item(9, 14).Value = 1.014178534644776

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023024301193869
$ws.Cells.Item(10, 4).Value = 1.028536733768854
$ws.Cells.Item(10, 5).Value = 1.046207874438275
$ws.Cells.Item(10, 6).Value = 1.050900080087144
$ws.Cells.Item(10, 9).Value = 1.031502013880895
$ws.Cells.Item(10, 10).Value = 1.029447001954771
$ws.Cells.Item(10, 11).Value = 1.03200921137064
$ws.Cells.Item(10, 12).Value = 1.049617393597486
$ws.Cells.Item(10, 13).Value = 1.054293273700589
$ws.Cells.Item(10, 14).Value = 1.013753978582449

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022242570224112
$ws.Cells.Item(11, 4).Value = 1.02795505568326
$ws.Cells.Item(11, 5).Value = 1.045417717131939
$ws.Cells.Item(11, 6).Value = 1.050074221421904
$ws.Cells.Item(11, 9).Value = 1.03134366603433
$ws.Cells.Item(11, 10).Value = 1.028903955712889
$ws.Cells.Item(11, 11).Value = 1.031552241669065
$ws.Cells.Item(11, 12).Value = 1.048950255202594
$ws.Cells.Item(11, 13).Value = 1.05358992184788
$ws.Cells.Item(11, 14).Value = 1.013569579605935

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.021952111474384
$ws.Cells.Item(12, 4).Value = 1.027738889441527
$ws.Cells.Item(12, 5).Value = 1.04512430667064
$ws.Cells.Item(12, 6).Value = 1.049767516839058
$ws.Cells.Item(12, 9).Value = 1.031284562157452
$ws.Cells.Item(12, 10).Value = 1.028702070378827
$ws.Cells.Item(12, 11).Value = 1.031382280825497
$ws.Cells.Item(12, 12).Value = 1.048702435488214
$ws.Cells.Item(12, 13).Value = 1.053328620032368
$ws.Cells.Item(12, 14).Value = 1.013501000626008

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022014419930789
$ws.Cells.Item(13, 4).Value = 1.027785262603421
$ws.Cells.Item(13, 5).Value = 1.045187240204456
$ws.Cells.Item(13, 6).Value = 1.049833303505642
$ws.Cells.Item(13, 9).Value = 1.031297253106941
$ws.Cells.Item(13, 10).Value = 1.028745383353395
$ws.Cells.Item(13, 11).Value = 1.031418748061752
$ws.Cells.Item(13, 12).Value = 1.048755594353264
$ws.Cells.Item(13, 13).Value = 1.053384672250011
$ws.Cells.Item(13, 14).Value = 1.013515714901536

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022218562643898
$ws.Cells.Item(14, 4).Value = 1.027937189463071
$ws.Cells.Item(14, 5).Value = 1.045393461923185
$ws.Cells.Item(14, 6).Value = 1.05004886799219
$ws.Cells.Item(14, 9).Value = 1.031338786334594
$ws.Cells.Item(14, 10).Value = 1.028887271360888
$ws.Cells.Item(14, 11).Value = 1.031538197185707
$ws.Cells.Item(14, 12).Value = 1.048929770658779
$ws.Cells.Item(14, 13).Value = 1.053568323498587
$ws.Cells.Item(14, 14).Value = 1.013563912581831

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022344329834042
$ws.Cells.Item(15, 4).Value = 1.028030782664014
$ws.Cells.Item(15, 5).Value = 1.045520533690917
$ws.Cells.Item(15, 6).Value = 1.05018169176265
$ws.Cells.Item(15, 9).Value = 1.031364338361246
$ws.Cells.Item(15, 10).Value = 1.028974670178186
$ws.Cells.Item(15, 11).Value = 1.031611764295484
$ws.Cells.Item(15, 12).Value = 1.049037084420324
$ws.Cells.Item(15, 13).Value = 1.053681470999048
$ws.Cells.Item(15, 14).Value = 1.013593597491052

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023076169745374
$ws.Cells.Item(16, 4).Value = 1.028575323133628
$ws.Cells.Item(16, 5).Value = 1.04626032717258
$ws.Cells.Item(16, 6).Value = 1.050954897579276
$ws.Cells.Item(16, 9).Value = 1.031512482755844
$ws.Cells.Item(16, 10).Value = 1.029483017847291
$ws.Cells.Item(16, 11).Value = 1.032039507869459
$ws.Cells.Item(16, 12).Value = 1.049661667383655
$ws.Cells.Item(16, 13).Value = 1.054339946540417
$ws.Cells.Item(16, 14).Value = 1.013766204635627

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023535077448298
$ws.Cells.Item(17, 4).Value = 1.028916712570257
$ws.Cells.Item(17, 5).Value = 1.046724540796878
$ws.Cells.Item(17, 6).Value = 1.051440011195907
$ws.Cells.Item(17, 9).Value = 1.031604899747806
$ws.Cells.Item(17, 10).Value = 1.029801582643508
$ws.Cells.Item(17, 11).Value = 1.03230742557338
$ws.Cells.Item(17, 12).Value = 1.050053427048592
$ws.Cells.Item(17, 13).Value = 1.054752911345847
$ws.Cells.Item(17, 14).Value = 1.013874325456059

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.023802695046705
$ws.Cells.Item(18, 4).Value = 1.029115772335578
$ws.Cells.Item(18, 5).Value = 1.046995367228865
$ws.Cells.Item(18, 6).Value = 1.051723007153716
$ws.Cells.Item(18, 9).Value = 1.03165862130137
$ws.Cells.Item(18, 10).Value = 1.029987285428291
$ws.Cells.Item(18, 11).Value = 1.032463555544461
$ws.Cells.Item(18, 12).Value = 1.050281925421228
$ws.Cells.Item(18, 13).Value = 1.054993758614751
$ws.Cells.Item(18, 14).Value = 1.013937336165306

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023893936389662
$ws.Cells.Item(19, 4).Value = 1.029183635245166
$ws.Cells.Item(19, 5).Value = 1.047087722014213
$ws.Cells.Item(19, 6).Value = 1.051819507862128
$ws.Cells.Item(19, 9).Value = 1.031676907811664
$ws.Cells.Item(19, 10).Value = 1.030050586546393
$ws.Cells.Item(19, 11).Value = 1.032516767789389
$ws.Cells.Item(19, 12).Value = 1.050359836080907
$ws.Cells.Item(19, 13).Value = 1.055075876590407
$ws.Cells.Item(19, 14).Value = 1.013958811982726

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023485846753951
$ws.Cells.Item(20, 4).Value = 1.02888009164286
$ws.Cells.Item(20, 5).Value = 1.046674729033092
$ws.Cells.Item(20, 6).Value = 1.051387959243247
$ws.Cells.Item(20, 9).Value = 1.0315950032844
$ws.Cells.Item(20, 10).Value = 1.029767415113314
$ws.Cells.Item(20, 11).Value = 1.032278695223238
$ws.Cells.Item(20, 12).Value = 1.05001139582633
$ws.Cells.Item(20, 13).Value = 1.054708607051629
$ws.Cells.Item(20, 14).Value = 1.013862730738287

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022158450180733
$ws.Cells.Item(21, 4).Value = 1.027892453697486
$ws.Cells.Item(21, 5).Value = 1.045332732320058
$ws.Cells.Item(21, 6).Value = 1.049985388090625
$ws.Cells.Item(21, 9).Value = 1.031326563743992
$ws.Cells.Item(21, 10).Value = 1.02884549368692
$ws.Cells.Item(21, 11).Value = 1.031503028530127
$ws.Cells.Item(21, 12).Value = 1.048878480499309
$ws.Cells.Item(21, 13).Value = 1.05351424405601
$ws.Cells.Item(21, 14).Value = 1.013549721915198

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021323346804629
$ws.Cells.Item(22, 4).Value = 1.027270878027396
$ws.Cells.Item(22, 5).Value = 1.044489480623814
$ws.Cells.Item(22, 6).Value = 1.049103859542643
$ws.Cells.Item(22, 9).Value = 1.031156127750871
$ws.Cells.Item(22, 10).Value = 1.028264838999029
$ws.Cells.Item(22, 11).Value = 1.031014051506109
$ws.Cells.Item(22, 12).Value = 1.048166085212567
$ws.Cells.Item(22, 13).Value = 1.05276303586247
$ws.Cells.Item(22, 14).Value = 1.013352428947739

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.02176610065918
$ws.Cells.Item(23, 4).Value = 1.027600445003567
$ws.Cells.Item(23, 5).Value = 1.044936456074923
$ws.Cells.Item(23, 6).Value = 1.049571144474615
$ws.Cells.Item(23, 9).Value = 1.031246636310851
$ws.Cells.Item(23, 10).Value = 1.028572750842028
$ws.Cells.Item(23, 11).Value = 1.031273389547785
$ws.Cells.Item(23, 12).Value = 1.048543748145652
$ws.Cells.Item(23, 13).Value = 1.053161291193807
$ws.Cells.Item(23, 14).Value = 1.013457064424941

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023508092157125
$ws.Cells.Item(24, 4).Value = 1.028896639270518
$ws.Cells.Item(24, 5).Value = 1.046697236644426
$ws.Cells.Item(24, 6).Value = 1.05141147916465
$ws.Cells.Item(24, 9).Value = 1.031599475637921
$ws.Cells.Item(24, 10).Value = 1.029782854291506
$ws.Cells.Item(24, 11).Value = 1.032291677670294
$ws.Cells.Item(24, 12).Value = 1.050030387952523
$ws.Cells.Item(24, 13).Value = 1.054728626342077
$ws.Cells.Item(24, 14).Value = 1.013867970060087

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025527407399447
$ws.Cells.Item(25, 4).Value = 1.030398178219482
$ws.Cells.Item(25, 5).Value = 1.048742824453866
$ws.Cells.Item(25, 6).Value = 1.053548558608132
$ws.Cells.Item(25, 9).Value = 1.032001720393711
$ws.Cells.Item(25, 10).Value = 1.031182775050066
$ws.Cells.Item(25, 11).Value = 1.033467781228611
$ws.Cells.Item(25, 12).Value = 1.051755226045945
$ws.Cells.Item(25, 13).Value = 1.056546334361985
$ws.Cells.Item(25, 14).Value = 1.014342672831431

